$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 9556.5
$ws.Range("I62").Value = 12458.75
$ws.Range("J62").Value = 3752
$ws.Range("K62").Value = 12458.75
$ws.Range("L62").Value = 3752
$ws.Range("M62").Value = -11834.75
$ws.Range("N62").Value = -5000
$ws.Range("H65").Value = 9556.5
$ws.Range("I65").Value = 12458.75
$ws.Range("J65").Value = 3752
$ws.Range("K65").Value = 62293.75
$ws.Range("L65").Value = 18760
$ws.Range("M65").Value = -59173.75
$ws.Range("N65").Value = -25000
$ws.Range("H74").Value = 1964334.9
$ws.Range("I74").Value = 2177425.2
$ws.Range("J74").Value = 3903.2
$ws.Range("K74").Value = 2177425.2
$ws.Range("L74").Value = 3903.2
$ws.Range("M74").Value = -2176489.2
$ws.Range("N74").Value = -5775.2
$ws.Range("H77").Value = 1964334.9
$ws.Range("I77").Value = 2177425.2
$ws.Range("J77").Value = 3903.2
$ws.Range("K77").Value = 10887126
$ws.Range("L77").Value = 19516
$ws.Range("M77").Value = -10882446
$ws.Range("N77").Value = -28876

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6183.523
$ws.Range("I32").Value = 4398.737
$ws.Range("K32").Value = 4398.737
$ws.Range("M32").Value = -4111.737
$ws.Range("H61").Value = 5500.25
$ws.Range("I61").Value = 1000.3333
$ws.Range("J61").Value = 19000
$ws.Range("K61").Value = 1000.3333
$ws.Range("L61").Value = 19000
$ws.Range("M61").Value = -788.3333
$ws.Range("N61").Value = -19424
$ws.Range("H136").Value = 5500.25
$ws.Range("I136").Value = 1000.3333
$ws.Range("J136").Value = 19000
$ws.Range("K136").Value = 3000.9999
$ws.Range("L136").Value = 57000
$ws.Range("M136").Value = -450.9998999999998
$ws.Range("N136").Value = -62100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 468529.66
$ws.Range("I86").Value = 2083.3333
$ws.Range("J86").Value = 779493.9
$ws.Range("K86").Value = 2083.3333
$ws.Range("L86").Value = 779493.9
$ws.Range("M86").Value = -960.3332999999998
$ws.Range("N86").Value = -781739.9
$ws.Range("H89").Value = 468529.66
$ws.Range("I89").Value = 2083.3333
$ws.Range("J89").Value = 779493.9
$ws.Range("K89").Value = 10416.6665
$ws.Range("L89").Value = 3897469.5
$ws.Range("M89").Value = -4800.666499999999
$ws.Range("N89").Value = -3908701.5
$ws.Range("H99").Value = 1300
$ws.Range("I99").Value = 1066.6666
$ws.Range("K99").Value = 1066.6666
$ws.Range("M99").Value = 431.3334
$ws.Range("H134").Value = 3230.96
$ws.Range("I134").Value = 1949.6923
$ws.Range("K134").Value = 5849.0769
$ws.Range("M134").Value = -3314.0769

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 842.0769
$ws.Range("I7").Value = 1502.5714
$ws.Range("J7").Value = 71.5
$ws.Range("K7").Value = 1502.5714
$ws.Range("L7").Value = 71.5
$ws.Range("M7").Value = -1389.5714
$ws.Range("N7").Value = -297.5
$ws.Range("H22").Value = 1488.7
$ws.Range("I22").Value = 2197.4
$ws.Range("J22").Value = 780
$ws.Range("K22").Value = 2197.4
$ws.Range("L22").Value = 780
$ws.Range("M22").Value = -1847.4
$ws.Range("N22").Value = -1480
$ws.Range("H31").Value = 22383.523
$ws.Range("I31").Value = 22659.309
$ws.Range("J31").Value = 20487.5
$ws.Range("K31").Value = 22659.309
$ws.Range("L31").Value = 20487.5
$ws.Range("M31").Value = -22364.309
$ws.Range("N31").Value = -21077.5
$ws.Range("H34").Value = 22383.523
$ws.Range("I34").Value = 22659.309
$ws.Range("J34").Value = 20487.5
$ws.Range("K34").Value = 22659.309
$ws.Range("L34").Value = 20487.5
$ws.Range("M34").Value = -22457.309
$ws.Range("N34").Value = -20891.5
$ws.Range("H62").Value = 166669000
$ws.Range("I62").Value = 250001500
$ws.Range("K62").Value = 250001500
$ws.Range("M62").Value = -250000876
$ws.Range("H65").Value = 166669000
$ws.Range("I65").Value = 250001500
$ws.Range("K65").Value = 1250007500
$ws.Range("M65").Value = -1250004380
$ws.Range("H74").Value = 11890.4
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 11890.4
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 11890.4
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -13638.4
$ws.Range("H77").Value = 11890.4
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 11890.4
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 35671.2
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -44407.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 37037170
$ws.Range("I34").Value = 94
$ws.Range("J34").Value = 43478400
$ws.Range("K34").Value = 282
$ws.Range("L34").Value = 130435200
$ws.Range("M34").Value = -198
$ws.Range("N34").Value = -130435368
$ws.Range("H103").Value = 2606.5557
$ws.Range("I103").Value = 1749.75
$ws.Range("J103").Value = 3292
$ws.Range("K103").Value = 5249.25
$ws.Range("L103").Value = 9876
$ws.Range("M103").Value = -4370.25
$ws.Range("N103").Value = -11634
$ws.Range("H129").Value = 57224.168
$ws.Range("I129").Value = 1000
$ws.Range("J129").Value = 60531.47
$ws.Range("K129").Value = 3000
$ws.Range("L129").Value = 181594.41
$ws.Range("M129").Value = 2000
$ws.Range("N129").Value = -191594.41

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1026
$ws.Range("I2").Value = 28.75
$ws.Range("K2").Value = 28.75
$ws.Range("M2").Value = 84.25
$ws.Range("H80").Value = 2879
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 2813.8462
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 2813.8462
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -4809.8462
$ws.Range("H83").Value = 2879
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 2813.8462
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 14069.231
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -24053.231

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 696.6667
$ws.Range("I46").Value = 300
$ws.Range("J46").Value = 895
$ws.Range("K46").Value = 300
$ws.Range("L46").Value = 895
$ws.Range("M46").Value = -112
$ws.Range("N46").Value = -1271
$ws.Range("H64").Value = 20000
$ws.Range("J64").Value = 20000
$ws.Range("L64").Value = 20000
$ws.Range("N64").Value = -20450
$ws.Range("H67").Value = 20000
$ws.Range("J67").Value = 20000
$ws.Range("L67").Value = 20000
$ws.Range("N67").Value = -21560
$ws.Range("H68").Value = 2152.2104
$ws.Range("I68").Value = 1799
$ws.Range("J68").Value = 2637.875
$ws.Range("K68").Value = 1799
$ws.Range("L68").Value = 2637.875
$ws.Range("M68").Value = -1050
$ws.Range("N68").Value = -4135.875
$ws.Range("H71").Value = 2152.2104
$ws.Range("I71").Value = 1799
$ws.Range("J71").Value = 2637.875
$ws.Range("K71").Value = 8995
$ws.Range("L71").Value = 13189.375
$ws.Range("M71").Value = -5251
$ws.Range("N71").Value = -20677.375
$ws.Range("H82").Value = 2927.4
$ws.Range("I82").Value = 3950
$ws.Range("J82").Value = 2671.75
$ws.Range("K82").Value = 3950
$ws.Range("L82").Value = 2671.75
$ws.Range("M82").Value = -3589
$ws.Range("N82").Value = -3393.75
$ws.Range("H85").Value = 2927.4
$ws.Range("I85").Value = 3950
$ws.Range("J85").Value = 2671.75
$ws.Range("K85").Value = 3950
$ws.Range("L85").Value = 2671.75
$ws.Range("M85").Value = -2702
$ws.Range("N85").Value = -5167.75
$ws.Range("H136").Value = 2292.4736
$ws.Range("I136").Value = 1688.3572
$ws.Range("J136").Value = 3984
$ws.Range("K136").Value = 5065.071599999999
$ws.Range("L136").Value = 11952
$ws.Range("M136").Value = -2515.071599999999
$ws.Range("N136").Value = -17052
